# Applies the shared-string / field-table corrections to the "classFields" sheet.
# The PmsBrandDto field rows had their Field Name (column B) values rotated relative
# to their Field Modifier/Field Type, and the DemoController LOGGER/demoService rows
# (and the UmsAdminLoginParam username/password rows) had their names swapped.
# This fixes the Field Name (and, where needed, Field Type) values so each row
# describes the correct field again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# --- PmsBrandDto fields (rows 2-9): fix Field Name, and Field Type where it changes ---
$ws.Range("B2").Value = "logo"

$ws.Range("B3").Value = "name"

$ws.Range("B4").Value = "firstLetter"
$ws.Range("D4").Value = "java.lang.String"

$ws.Range("B5").Value = "brandStory"

$ws.Range("B6").Value = "showStatus"

$ws.Range("B7").Value = "sort"
$ws.Range("D7").Value = "java.lang.Integer"

$ws.Range("B8").Value = "factoryStatus"
$ws.Range("D8").Value = "java.lang.Integer"

$ws.Range("B9").Value = "bigPic"
$ws.Range("D9").Value = "java.lang.String"

# --- DemoController fields (rows 11-12): swap LOGGER / demoService ---
$ws.Range("B11").Value = "demoService"
$ws.Range("D11").Value = "com.macro.mall.demo.service.DemoService"

$ws.Range("B12").Value = "LOGGER"
$ws.Range("D12").Value = "org.slf4j.Logger"

# --- UmsAdminLoginParam fields (rows 13-14): swap password / username ---
$ws.Range("B13").Value = "username"

$ws.Range("B14").Value = "password"
